$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("spawning")

# New "site_code" column (D) on the spawning-sites table (rows 15-25).
$ws2.Range("D15").Value = "site_code"
$ws2.Range("D16").Value = "NSW02"
$ws2.Range("D17").Value = "NSW05"
$ws2.Range("D18").Value = "NSW09"
$ws2.Range("D19").Value = "NSW18"
$ws2.Range("D20").Value = "NSW13"
$ws2.Range("D21").Value = "NSW19"
$ws2.Range("D22").Value = "NSW19"
$ws2.Range("D23").Value = "NSW07"
$ws2.Range("D24").Value = "NSW10"
$ws2.Range("D25").Value = "NSW11"

# Move the active tab / selection from "basic-biology" to "spawning".
$ws2.Activate()
[void]$ws2.Range("P20").Select()
